# "add almost all lineal" -- refresh the generated numeric results across the
# linear-experiment workbook (new Lambda/Beta/Gamma table, new modified point,
# new vec_bf / vec_BF vectors). All of these cells hold their numbers as plain
# text (shared strings), so we force text formatting before writing the
# numeric-looking values to keep Excel from "helpfully" re-typing them as
# real numbers. NumberFormat is applied per contiguous single-area range
# only (multi-area comma ranges only honour their first area here).
#
# Sheets are addressed by index rather than by name: "Vector_bf" and
# "Vector_BF" differ only by case and Worksheets.Item(name) resolves
# case-insensitively to whichever of the two comes first, so a name lookup
# for "Vector_BF" would silently hit the "Vector_bf" sheet instead.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Restricciones_del_follower (sheet 3): refresh Lambda/Beta/Gamma table
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

# Column B (Function_Evaluation), rows 2-6 all change.
$ws3.Range("B2:B6").NumberFormat = "@"
$ws3.Range("B2").Value = "-4.382729079133727"
$ws3.Range("B3").Value = "0.38272907913372656"
$ws3.Range("B4").Value = "-13.81788203899294"
$ws3.Range("B5").Value = "3.3894519012124054"
$ws3.Range("B6").Value = "-8.517138863787123"

# Column D (Lambda_value), rows 2-6 all change.
$ws3.Range("D2:D6").NumberFormat = "@"
$ws3.Range("D2").Value = "0.0866877650392671"
$ws3.Range("D3").Value = "0.9648587319705634"
$ws3.Range("D4").Value = "0.9761226555169311"
$ws3.Range("D5").Value = "0.8143958706897286"
$ws3.Range("D6").Value = "0.23927405565041526"

# Column E (Beta_value), rows 2-6 all change.
$ws3.Range("E2:E6").NumberFormat = "@"
$ws3.Range("E2").Value = "0.746488456643626"
$ws3.Range("E3").Value = "0"
$ws3.Range("E4").Value = "0"
$ws3.Range("E5").Value = "0.09323959926382419"
$ws3.Range("E6").Value = "-0.6770049979448282"

# Column F (Gamma_value): rows 2,3,4,6 change; row 5 stays "0" untouched.
$ws3.Range("F2:F4").NumberFormat = "@"
$ws3.Range("F2").Value = "0"
$ws3.Range("F3").Value = "0.6119990159669402"
$ws3.Range("F4").Value = "0.607732950283753"
$ws3.Range("F6").NumberFormat = "@"
$ws3.Range("F6").Value = "-0.5008550602329267"

# Column A (Expression): only row 4's expression text changes.
$ws3.Range("A4").Value = "-16 - 2x + y_1 + 4y_2"

# ---------------------------------------------------------------------
# Punto_modificado (sheet 4): refresh the modified (x, y_1, y_2) point
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A2:C2").NumberFormat = "@"
$ws4.Range("A2").Value = "5.875840352759835"
$ws4.Range("B2").Value = "4.382729079133727"
$ws4.Range("C2").Value = "2.387767396848251"

# ---------------------------------------------------------------------
# Vector_bf (sheet 5): refresh vec_bf
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("A2:A3").NumberFormat = "@"
$ws5.Range("A2").Value = "-1.1901413818371251"
$ws5.Range("A3").Value = "-3.9044906220677245"

# ---------------------------------------------------------------------
# Vector_BF (sheet 6): refresh vec_BF (3rd entry -2.0 is unchanged)
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)
$ws6.Range("A2:A3").NumberFormat = "@"
$ws6.Range("A2").Value = "-1.09992679000025"
$ws6.Range("A3").Value = "2.2992388614901453"
